# Apply numeric value updates to Ultros_Profits workbook sheets
# (recomputed profit snapshot values in columns H-N across several sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 953.8
$ws.Range("J17").Value = 879.7368
$ws.Range("L17").Value = 2639.2104
$ws.Range("N17").Value = -2975.2104
$ws.Range("H33").Value = 569.46155
$ws.Range("I33").Value = 338.625
$ws.Range("J33").Value = 938.8
$ws.Range("K33").Value = 338.625
$ws.Range("L33").Value = 938.8
$ws.Range("M33").Value = -109.625
$ws.Range("N33").Value = -1396.8
$ws.Range("H48").Value = 805
$ws.Range("J48").Value = 665
$ws.Range("L48").Value = 1995
$ws.Range("N48").Value = -2579
$ws.Range("H56").Value = 805
$ws.Range("J56").Value = 665
$ws.Range("L56").Value = 1995
$ws.Range("N56").Value = -3063
$ws.Range("H58").Value = 428.75
$ws.Range("I58").Value = 71.666664
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 214.999992
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -64.99999199999999
$ws.Range("N58").Value = -4800
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H127").Value = 13382.4
$ws.Range("I127").Value = 1261.2858
$ws.Range("J127").Value = 41665
$ws.Range("K127").Value = 3783.8574
$ws.Range("L127").Value = 124995
$ws.Range("M127").Value = 1176.1426
$ws.Range("N127").Value = -134915

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18871340
$ws.Range("I32").Value = 20003412
$ws.Range("K32").Value = 20003412
$ws.Range("M32").Value = -20003125
$ws.Range("H36").Value = 252356.5
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3100.913
$ws.Range("I86").Value = 2093.75
$ws.Range("K86").Value = 2093.75
$ws.Range("M86").Value = -970.75
$ws.Range("H89").Value = 3100.913
$ws.Range("I89").Value = 2093.75
$ws.Range("K89").Value = 10468.75
$ws.Range("M89").Value = -4852.75
$ws.Range("H94").Value = 3058.261
$ws.Range("I94").Value = 2755.4666
$ws.Range("J94").Value = 3626
$ws.Range("K94").Value = 2755.4666
$ws.Range("L94").Value = 3626
$ws.Range("M94").Value = -2304.4666
$ws.Range("N94").Value = -4528
$ws.Range("H134").Value = 2096.8484
$ws.Range("I134").Value = 1639.8667
$ws.Range("K134").Value = 4919.6001
$ws.Range("M134").Value = -2384.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 41.6
$ws.Range("I7").Value = 35.11111
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 35.11111
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = 77.88889
$ws.Range("N7").Value = -326

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1097.6
$ws.Range("I22").Value = 775.1111
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 2325.3333
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = -2156.3333
$ws.Range("N22").Value = -12338
$ws.Range("H27").Value = 1097.6
$ws.Range("I27").Value = 775.1111
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 2325.3333
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = -2223.3333
$ws.Range("N27").Value = -12204
$ws.Range("H34").Value = 994.3333
$ws.Range("J34").Value = 1499.6666
$ws.Range("L34").Value = 4498.9998
$ws.Range("N34").Value = -4666.9998
$ws.Range("H44").Value = 190.4
$ws.Range("I44").Value = 150.5
$ws.Range("K44").Value = 451.5
$ws.Range("M44").Value = -53.5
$ws.Range("H122").Value = 4308.9
$ws.Range("I122").Value = 701.7143
$ws.Range("J122").Value = 6251.231
$ws.Range("K122").Value = 6315.428699999999
$ws.Range("L122").Value = 56261.079
$ws.Range("M122").Value = -3865.428699999999
$ws.Range("N122").Value = -61161.079
$ws.Range("H129").Value = 1118.9166
$ws.Range("I129").Value = 653.375
$ws.Range("J129").Value = 2050
$ws.Range("K129").Value = 1960.125
$ws.Range("L129").Value = 6150
$ws.Range("M129").Value = 3039.875
$ws.Range("N129").Value = -16150
$ws.Range("H131").Value = 3843.3125
$ws.Range("I131").Value = 2303.5715
$ws.Range("K131").Value = 6910.7145
$ws.Range("M131").Value = -1870.7145
$ws.Range("H134").Value = 3355.4614
$ws.Range("I134").Value = 2051.75
$ws.Range("J134").Value = 19000
$ws.Range("K134").Value = 6155.25
$ws.Range("L134").Value = 57000
$ws.Range("M134").Value = -1085.25
$ws.Range("N134").Value = -67140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 15
$ws.Range("I22").Value = 15
$ws.Range("K22").Value = 15
$ws.Range("M22").Value = 514
$ws.Range("H33").Value = 25499
$ws.Range("J33").Value = 25499
$ws.Range("L33").Value = 25499
$ws.Range("N33").Value = -26003
$ws.Range("H36").Value = 3092.6667
$ws.Range("I36").Value = 711.2
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 711.2
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = -226.2
$ws.Range("N36").Value = -15970
$ws.Range("H97").Value = 11901.889
$ws.Range("I97").Value = 833.4
$ws.Range("J97").Value = 25737.5
$ws.Range("K97").Value = 833.4
$ws.Range("L97").Value = 25737.5
$ws.Range("M97").Value = -337.4
$ws.Range("N97").Value = -26729.5
$ws.Range("H102").Value = 3194.6365
$ws.Range("I102").Value = 2229.4644
$ws.Range("J102").Value = 8599.6
$ws.Range("K102").Value = 2229.4644
$ws.Range("L102").Value = 8599.6
$ws.Range("M102").Value = -607.4643999999998
$ws.Range("N102").Value = -11843.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 699
$ws.Range("I16").Value = 651.55554
$ws.Range("J16").Value = 805.75
$ws.Range("K16").Value = 651.55554
$ws.Range("L16").Value = 805.75
$ws.Range("M16").Value = -481.55554
$ws.Range("N16").Value = -1145.75
$ws.Range("H46").Value = 3899.6667
$ws.Range("J46").Value = 3899.6667
$ws.Range("L46").Value = 3899.6667
$ws.Range("N46").Value = -4275.6667
$ws.Range("H93").Value = 7411.4
$ws.Range("I93").Value = 5630.3335
$ws.Range("J93").Value = 10083
$ws.Range("K93").Value = 5630.3335
$ws.Range("L93").Value = 10083
$ws.Range("M93").Value = -4382.3335
$ws.Range("N93").Value = -12579
$ws.Range("H136").Value = 2103.7354
$ws.Range("I136").Value = 2061.4243
$ws.Range("K136").Value = 6184.2729
$ws.Range("M136").Value = -3634.2729

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13774.9
$ws.Range("I45").Value = 7968
$ws.Range("J45").Value = 15226.625
$ws.Range("K45").Value = 7968
$ws.Range("L45").Value = 15226.625
$ws.Range("M45").Value = -7477
$ws.Range("N45").Value = -16208.625
$ws.Range("H62").Value = 391.66666
$ws.Range("I62").Value = 337.5
$ws.Range("K62").Value = 337.5
$ws.Range("M62").Value = 286.5
$ws.Range("H65").Value = 391.66666
$ws.Range("I65").Value = 337.5
$ws.Range("K65").Value = 1687.5
$ws.Range("M65").Value = 1432.5
$ws.Range("H81").Value = 6071.143
$ws.Range("J81").Value = 6071.143
$ws.Range("L81").Value = 12142.286
$ws.Range("N81").Value = -14264.286
$ws.Range("H84").Value = 6071.143
$ws.Range("J84").Value = 6071.143
$ws.Range("L84").Value = 60711.43
$ws.Range("N84").Value = -71319.42999999999

